$wb = $excel.ActiveWorkbook
$wsResources = $wb.Worksheets.Item("Resources")

# The "TwoLetterISOLanguageName" column header becomes "Language"
# (resource items will now be written straight to .resx files).
$wsResources.Range("B1").Value = "Language"

# Resources is now the active/selected sheet, with B1 selected.
$wsResources.Activate()
$wsResources.Range("B1").Select()
